# "Minor change in figure"
#
# The deck originally had two slides:
#   1) a blank title slide (ctrTitle/subTitle placeholders only)
#   2) the terminology figure
#
# The edit removes the superfluous blank title slide, leaving the
# terminology figure as the sole, first slide of the presentation.

$p = $ppt.ActivePresentation
$p.Slides.Item(1).Delete()
